$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "book_path_1"
$ws.Range("B1").Value = "book_path_2"
$ws.Range("C1").Value = "main_title"
$ws.Range("D1").Value = "type_of_book"
$ws.Range("E1").Value = "class_info"
$ws.Range("F1").Value = "age_or_class"
$ws.Range("G1").Value = "output_path"
$ws.Range("H1").Value = "template"
$ws.Range("I1").Value = "quantity"

# Row 2
$ws.Range("A2").Value = "/Users/rrkhikmatullin/Desktop/Снимок экрана 2024-04-06 в 16.42.33.png"
$ws.Range("B2").Value = "/Users/rrkhikmatullin/Desktop/Снимок экрана 2024-04-06 в 16.42.33.png"
$ws.Range("C2").Value = "Я считаю ДО десяти"
$ws.Range("D2").Value = "Для детей 5-6 лет"
$ws.Range("E2").Value = "5–6"
$ws.Range("F2").Value = "лет"
$ws.Range("G2").Value = "ya_shitau_do_2"
$ws.Range("H2").Value = 2
$ws.Range("I2").Value = 10

# Row 3
$ws.Range("A3").Value = "/Users/rrkhikmatullin/Desktop/Снимок экрана 2024-04-06 в 16.42.13.png"
$ws.Range("B3").Value = "/Users/rrkhikmatullin/Desktop/Снимок экрана 2024-04-06 в 16.42.13.png"
$ws.Range("C3").Value = "Гимназия для дошколят"
$ws.Range("D3").Value = "Для детей 5-6 лет"
$ws.Range("E3").Value = "6–6"
$ws.Range("F3").Value = "лет"
$ws.Range("G3").Value = "гимназя_2"
$ws.Range("H3").Value = 2
$ws.Range("I3").Value = 10

# Row 4
$ws.Range("A4").Value = "/Users/rrkhikmatullin/Desktop/Снимок экрана 2024-04-06 в 16.41.46.png"
$ws.Range("B4").Value = "/Users/rrkhikmatullin/Desktop/Снимок экрана 2024-04-06 в 16.41.46.png"
$ws.Range("C4").Value = "От звука к слову"
$ws.Range("D4").Value = "Для детей 5-6 лет"
$ws.Range("E4").Value = "7–6"
$ws.Range("F4").Value = "лет"
$ws.Range("G4").Value = "от_слова_2"
$ws.Range("H4").Value = 2
$ws.Range("I4").Value = 10

# Row 5
$ws.Range("A5").Value = "/Users/rrkhikmatullin/Desktop/Снимок экрана 2024-04-06 в 16.42.33.png"
$ws.Range("C5").Value = "Я считаю ДО десяти"
$ws.Range("D5").Value = "Для детей 5-6 лет"
$ws.Range("E5").Value = "5–6"
$ws.Range("F5").Value = "лет"
$ws.Range("G5").Value = "ya_shitau_do"
$ws.Range("H5").Value = 1

# Row 6
$ws.Range("A6").Value = "/Users/rrkhikmatullin/Desktop/Снимок экрана 2024-04-06 в 16.42.13.png"
$ws.Range("C6").Value = "Гимназия для дошколят"
$ws.Range("D6").Value = "Для детей 5-6 лет"
$ws.Range("E6").Value = "6–6"
$ws.Range("F6").Value = "лет"
$ws.Range("G6").Value = "гимназя"
$ws.Range("H6").Value = 1

# Row 7
$ws.Range("A7").Value = "/Users/rrkhikmatullin/Desktop/Снимок экрана 2024-04-06 в 16.41.46.png"
$ws.Range("C7").Value = "От звука к слову"
$ws.Range("D7").Value = "Для детей 5-6 лет"
$ws.Range("E7").Value = "7–6"
$ws.Range("F7").Value = "лет"
$ws.Range("G7").Value = "от_слова"
$ws.Range("H7").Value = 1

# clear the old class_info number format (was on D2:D4) since that column is now type_of_book
$ws.Range("D2:D7").ClearFormats()

# style E2:E7 with the numFmtId 16 "d-mmm" style used previously on the class_info column
$ws.Range("E2:E7").NumberFormat = "d-mmm"

$ws.Range("D11").Select()
